# "Generate Report for Handoff" — b.md has now been handed off for
# localization, so its status moves from "Handed back: in sync with en-US"
# to "Ready for handoff", a new handoff xliff file is recorded (with its
# generation datetime), and an error detail explaining the stale handback
# version is attached.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$overviewDate = "2016-08-31 18:44:30"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60d969a802cb1ad36d8afb8bc982bba1178b1c81/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba8b85c71f6d9e882005488793105af28bd7879a/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet — row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $overviewDate

# ---------------------------------------------------------------------
# zh-cn sheet — row 3 is "b.md".
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 18:44:25"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet — row 3 is "b.md".
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
